$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.368.75"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "3.348.05"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "557.20"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "174.25"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").Value = "3.334.10"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +6.68%  "
$ws.Range("D11").Value = "0.635"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").Value = "53.96"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "0.0000276"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "9.09"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "3.877.53"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "18.20"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "3.346.03"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "64.883.77"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "11.74"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "0.989"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "454.94"
$ws.Range("E22").Value = "  +5.53%  "
$ws.Range("D23").Value = "4.93"
$ws.Range("E23").Value = "  +7.84%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "14.12"
$ws.Range("E25").Value = "  +8.09%  "
$ws.Range("D26").Value = "87.07"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").Value = "10.79"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "31.09"
$ws.Range("E30").Value = "  +5.56%  "
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").Value = "11.47"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "568.51"
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("D34").Value = "61.10"
$ws.Range("E34").Value = "  +4.47%  "
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "3.61"
$ws.Range("E37").Value = "  +4.24%  "
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "35.43"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.369"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").Value = "3.063.96"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "0.0416"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").Value = "3.21"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "0.134"
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "142.15"
$ws.Range("E49").Value = "  +6.67%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").Value = "8.12"
$ws.Range("E51").Value = "  -0.43%  "
